# Slide 2 ("Introductions"), body placeholder shape (idx 4 in Shapes),
# paragraph 4: "SDLC experience" -> split into three runs:
#   "SDLC" + "/DevOps " + "experience"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(4)

# Select the single space character between "SDLC" and "experience" and
# insert "/DevOps" immediately before it. This splits the original single
# run into three runs: "SDLC", "/DevOps " (the inserted text plus the
# pre-existing space), and "experience" - matching the target XML.
$spaceChar = $para.Characters(5, 1)
$spaceChar.InsertBefore("/DevOps") | Out-Null
